$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPG")

# Row 4 (Inventory)
$ws.Range("B4").Value = 1914000000.0
$ws.Range("C4").Value = 1735000000.0
$ws.Range("D4").Value = 1672000000.0
$ws.Range("E4").Value = 1706000000.0
$ws.Range("F4").Value = 1859000000.0

# Row 14 (Accounts Payable)
$ws.Range("B14").Value = 3714000000.0
$ws.Range("C14").Value = 2259000000.0
$ws.Range("D14").Value = 3401000000.0
$ws.Range("E14").Value = 2970000000.0
$ws.Range("F14").Value = 3068000000.0

# Row 20 (Long Term Tax Liability (Deferred))
$ws.Range("B20").Value = 88000000.0
$ws.Range("C20").Value = 56000000.0
$ws.Range("D20").Value = 179000000.0
$ws.Range("E20").Value = 209000000.0
$ws.Range("F20").Value = 216000000.0
